# Abstract template update: drop the "and" before Chlorophyll-a and append a
# new clause about nutrients / other parameters, as a separate run (matching
# the author's "Finished update to discretewq EDI publication" edit).

$d = $word.ActiveDocument

# Step 1: " (surface), and Chlorophyll-a concentration (surface)"
#      -> " (surface), Chlorophyll-a concentration (surface)"
# Remove the now-superfluous "and " right before "Chlorophyll-a".
$findAnd = $d.Content
$findAnd.Find.Execute("and Chlorophyll-a", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$andRange = $d.Range($findAnd.Start, $findAnd.Start + 4)
$andRange.Text = ""

# Step 2: append ", nutrients (surface), and other parameters" right after
# "...Chlorophyll-a concentration (surface)" as its own run (toggling Bold
# on/off forces the new text to stay in a distinct <w:r> instead of being
# silently re-absorbed into the neighbouring run).
$findConc = $d.Content
$findConc.Find.Execute("Chlorophyll-a concentration (surface)", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)
$findConc.Collapse(0)
$findConc.InsertAfter(", nutrients (surface), and other parameters")
$findConc.Bold = 1
$findConc.Bold = 0
